# Weekly driver report update for 2025-04-19
# Refreshes the "Good Drivers (Roaming > 99.8%)" table (rows 12-54) on the
# "Driver Summary" sheet with this week's data. One new row (54) appears,
# and the table now spans one row further than before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(12, 'Intel(R) Wireless-AC 9560 160MHz - 22.220.0.4', 226852, 99.90000000000001, $null)
    ,@(13, 'Intel(R) Wireless-AC 9560 160MHz - 23.40.1.1', 276085, 99.90000000000001, $null)
    ,@(14, 'Intel(R) Wireless-AC 9560 160MHz - 23.30.0.6', 625139, 99.90000000000001, $null)
    ,@(15, 'Intel(R) Wireless-AC 9560 160MHz - 23.90.0.2', 4931894, 99.90000000000001, $null)
    ,@(16, 'Intel(R) Wireless-AC 9560 160MHz - 22.190.0.4', 27599, 100, $null)
    ,@(17, 'Intel(R) Wireless-AC 9560 160MHz - 21.0.1.1', 191877, 100, $null)
    ,@(18, 'Intel(R) Wireless-AC 9560 160MHz - 21.90.1.2', 262523, 100, $null)
    ,@(19, 'Intel(R) Wireless-AC 9560 160MHz - 21.80.2.3', 10451, 100, $null)
    ,@(20, 'Intel(R) Wireless-AC 9560 160MHz - 21.10.2.2', 61902, 100, $null)
    ,@(21, 'Intel(R) Wireless-AC 9560 160MHz - 21.110.1.1', 684542, 99.90000000000001, '2025-02-05')
    ,@(22, 'Intel(R) Wireless-AC 9560 160MHz - 23.110.0.5', 53193, 100, '2025-01-01')
    ,@(23, 'Intel(R) Wireless-AC 9560 160MHz - 23.50.0.6', 27295, 100, '2025-01-01')
    ,@(24, 'Intel(R) Wireless-AC 9560 160MHz - 20.110.0.3', 59635, 99.90000000000001, '2024-04-13')
    ,@(25, 'Intel(R) Wireless-AC 9560 160MHz - 23.70.4.1', 52237, 100, '2024-04-13')
    ,@(26, 'Intel(R) Wireless-AC 9560 160MHz - 20.100.0.4', 108823, 100, '2023-12-19')
    ,@(27, 'Intel(R) Wireless-AC 9560 160MHz - 22.250.10.1', 78331, 99.90000000000001, '2023-08-14')
    ,@(28, 'Intel(R) Wireless-AC 9560 160MHz - 22.130.0.5', 109036, 99.90000000000001, '2023-03-28')
    ,@(29, 'Intel(R) Wireless-AC 9560 160MHz - 22.150.1.1', 154405, 100, '2022-08-29')
    ,@(30, 'Intel(R) Wireless-AC 9560 160MHz - 22.160.0.4', 3650830, 99.90000000000001, '2022-08-13')
    ,@(31, 'Intel(R) Wireless-AC 9560 160MHz - 22.120.0.3', 99547, 99.90000000000001, '2022-01-30')
    ,@(32, 'Intel(R) Wireless-AC 9560 160MHz - 23.20.0.4', 44160, 99.90000000000001, '2021-08-18')
    ,@(33, 'Intel(R) Wireless-AC 9560 160MHz - 22.70.0.6', 75637, 100, '2021-06-28')
    ,@(34, 'Intel(R) Wireless-AC 9560 160MHz - 22.80.1.1', 93224, 100, '2021-06-28')
    ,@(35, 'Intel(R) Wireless-AC 9560 160MHz - 22.200.2.1', 450036, 100, '2021-04-27')
    ,@(36, 'Intel(R) Wireless-AC 9560 160MHz - 22.50.0.7', 1543020, 100, '2021-04-18')
    ,@(37, 'Intel(R) Wireless-AC 9560 160MHz - 22.40.0.7', 169789, 100, '2021-03-02')
    ,@(38, 'Intel(R) Wireless-AC 9560 160MHz - 22.30.0.11', 238746, 99.90000000000001, '2021-01-19')
    ,@(39, 'Intel(R) Wireless-AC 9560 160MHz - 22.10.0.7', 321983, 99.90000000000001, '2020-10-19')
    ,@(40, 'Intel(R) Wireless-AC 9560 160MHz - 21.120.0.9', 95547, 100, '2020-08-15')
    ,@(41, 'Intel(R) Wireless-AC 9560 160MHz - 21.30.3.2', 443012, 100, '2020-08-15')
    ,@(42, 'Intel(R) Wireless-AC 9560 160MHz - 21.110.3.2', 121232, 99.90000000000001, '2020-08-05')
    ,@(43, 'Intel(R) Wireless-AC 9560 160MHz - 21.110.2.1', 36791, 99.90000000000001, '2020-07-20')
    ,@(44, 'Intel(R) Wireless-AC 9560 160MHz - 22.100.0.3', 25808, 100, '2020-04-05')
    ,@(45, 'Intel(R) Wireless-AC 9560 160MHz - 21.80.2.1', 209593, 99.90000000000001, '2020-02-24')
    ,@(46, 'Intel(R) Wireless-AC 9560 160MHz - 22.140.0.3', 178916, 99.90000000000001, '2020-02-24')
    ,@(47, 'Intel(R) Wireless-AC 9560 160MHz - 21.70.0.6', 67365, 100, '2020-01-06')
    ,@(48, 'Intel(R) Wireless-AC 9560 160MHz - 23.10.0.8', 97122, 100, '2020-01-06')
    ,@(49, 'Intel(R) Wireless-AC 9560 160MHz - 21.50.1.1', 308468, 100, '2019-10-05')
    ,@(50, 'Intel(R) Wireless-AC 9560 160MHz - 21.40.2.2', 139827, 100, '2019-08-31')
    ,@(51, 'Intel(R) Wireless-AC 9560 160MHz - 21.0.0.5', 194659, 100, '2019-08-31')
    ,@(52, 'Intel(R) Wireless-AC 9560 160MHz - 21.40.1.4', 34718, 100, '2019-08-14')
    ,@(53, 'Intel(R) Wireless-AC 9560 160MHz - 21.10.1.2', 546751, 100, '2019-08-10')
    ,@(54, 'Intel(R) Wireless-AC 9560 160MHz - 21.10.0.5', 108633, 100, '2019-04-06')
)

foreach ($row in $data) {
    $r = $row[0]
    $driver = $row[1]
    $samples = $row[2]
    $pct = $row[3]
    $vintage = $row[4]

    $ws.Cells.Item($r, 1).Value = $driver
    $ws.Cells.Item($r, 2).Value = $samples
    $ws.Cells.Item($r, 4).Value = $pct

    $eCell = $ws.Cells.Item($r, 5)
    if ($null -eq $vintage) {
        $eCell.Value = ""
    } else {
        # Force text so the date-formatted string isn't reinterpreted as a
        # real Excel date serial number.
        $eCell.NumberFormat = "@"
        $eCell.Value = $vintage
        $eCell.NumberFormat = "General"
    }
}
